# Add a new default user "Oleg_Babak" to the Users sheet of CobaltUsers.xlsx
# New row is inserted before the "danholland" row (row 31), pushing subsequent
# rows down by one, matching the behavior of selecting row 31 and choosing
# Insert > Sheet Rows in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Insert a new blank row at row 31, shifting rows 31..end down by one.
$ws.Rows.Item(31).Insert()

# Match the thin-border "data row" formatting used throughout the table
# (Excel style index 5: thin border on all sides, no explicit fill).
$newRow = $ws.Range("A31:G31")
$newRow.Borders.LineStyle = 1   # xlContinuous
$newRow.Borders.Weight = 2      # xlThin
$newRow.Borders.Color = 0       # black

# Fill in the data for the newly inserted row.
$ws.Cells.Item(31, 1).Value = "Oleg_Babak"
$ws.Cells.Item(31, 2).Value = "Password1!"
$ws.Cells.Item(31, 4).Value = "CUSTOM_USER"
$ws.Cells.Item(31, 5).Value = "Smoke Test User"
$ws.Cells.Item(31, 6).Value = "N"

# Update the visible selection/scroll position to match the saved workbook state.
$ws.Range("E30").Select()
$ws.Application.ActiveWindow.ScrollRow = 12
